$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.432.73'
$ws.Range("D3").Value = '1.848.63'
$ws.Range("E3").Value = '  +0.34%  '
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").Value = "'240.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.86%  '
$ws.Range("D6").Value = "'0.6297"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.11%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("D8").Value = "'0.07677"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.02%  '
$ws.Range("D9").Value = "'0.2926"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.57%  '
$ws.Range("E10").Value = '  +0.93%  '
$ws.Range("D11").Value = "'0.07737"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.55%  '
$ws.Range("D12").Value = '1.862.28'
$ws.Range("E12").Value = '  +1.56%  '
$ws.Range("D14").Value = "'0.6798"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.21%  '
$ws.Range("E15").Value = '  +2.66%  '
$ws.Range("D16").Value = "'83.60"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.74%  '
$ws.Range("D17").Value = '2.118.90'
$ws.Range("E17").Value = '  +1.42%  '
$ws.Range("D18").Value = "'6.192"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.82%  '
$ws.Range("D19").Value = '29.478.42'
$ws.Range("E19").Value = '  +0.38%  '
$ws.Range("D20").Value = "'228.19"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.05%  '
$ws.Range("E21").Value = '  +0.05%  '
$ws.Range("E22").Value = '  +0.05%  '
$ws.Range("D23").Value = "'7.430"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.05%  '
$ws.Range("E24").Value = '  +0.01%  '
$ws.Range("D25").Value = "'157.80"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.83%  '
$ws.Range("D26").Value = "'0.1379"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.86%  '
$ws.Range("D27").Value = "'8.405"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.77%  '
$ws.Range("E28").Value = '  +0.52%  '
$ws.Range("D29").Value = "'1.350"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +6.26%  '
$ws.Range("D30").Value = "'1.463"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.54%  '
$ws.Range("D31").Value = "'0.05670"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.69%  '
$ws.Range("D32").Value = "'4.121"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.45%  '
$ws.Range("D33").Value = "'4.029"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.22%  '
$ws.Range("D34").Value = "'1.842"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.70%  '
$ws.Range("D35").Value = "'1.161"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.55%  '
$ws.Range("D36").Value = "'0.7079"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.12%  '
$ws.Range("D37").Value = "'2.586"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Value = "'2.779"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.89%  '
$ws.Range("D40").Value = '1.219.34'
$ws.Range("E40").Value = '  -2.09%  '
$ws.Range("D41").Value = "'6.545"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.90%  '
$ws.Range("D42").Value = "'0.9110"
$ws.Range("D42").Style = "Normal"
$ws.Range("E43").Value = '  +0.15%  '
$ws.Range("D44").Value = "'101.66"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.16%  '
$ws.Range("D45").Value = "'65.98"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.50%  '
$ws.Range("D46").Value = "'0.00000000120"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.04%  '
$ws.Range("D47").Value = "'7.138"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.41%  '
$ws.Range("E48").Value = '  +0.76%  '
$ws.Range("D49").Value = "'9.042"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.62%  '
$ws.Range("D50").Value = "'1.678"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.33%  '
$ws.Range("D51").Value = "'0.1146"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.29%  '
